$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 280996
$ws.Range("E2").Value = 11884
$ws.Range("F2").Value = 11884
$ws.Range("G2").Value = 10780
$ws.Range("H2").Value = 6157
$ws.Range("I2").Value = 5266
$ws.Range("J2").Value = 891
$ws.Range("K2").Value = 400723
$ws.Range("L2").Value = 225284
$ws.Range("M2").Value = 175439
$ws.Range("N2").Value = 167258
$ws.Range("O2").Value = 8181
$ws.Range("P2").Value = 1575
$ws.Range("Q2").Value = 16299
$ws.Range("R2").Value = -6000
$ws.Range("S2").Value = -4162
$ws.Range("T2").Value = 18631
$ws.Range("U2").Value = -2333
$ws.Range("V2").Value = 127531
$ws.Range("W2").Value = 4.23
$ws.Range("X2").Value = 2.19
$ws.Range("Y2").Value = 3.2
$ws.Range("Z2").Value = 1.56
$ws.Range("AA2").Value = 128.41
$ws.Range("AB2").Value = 10333.68
$ws.Range("AC2").Value = 16724
$ws.Range("AD2").Value = 15.43
$ws.Range("AE2").Value = 565975
$ws.Range("AF2").Value = 0.46
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 0.78
$ws.Range("AI2").Value = 11.22
$ws.Range("AJ2").Value = 31490892

# Row 3
$ws.Range("D3").Value = 291277
$ws.Range("E3").Value = 8537
$ws.Range("F3").Value = 8537
$ws.Range("G3").Value = -798
$ws.Range("H3").Value = -3455
$ws.Range("I3").Value = -3831
$ws.Range("J3").Value = 376
$ws.Range("K3").Value = 406934
$ws.Range("L3").Value = 236055
$ws.Range("M3").Value = 170878
$ws.Range("N3").Value = 162475
$ws.Range("O3").Value = 8404
$ws.Range("P3").Value = 1575
$ws.Range("Q3").Value = 6411
$ws.Range("R3").Value = -15584
$ws.Range("S3").Value = 7310
$ws.Range("T3").Value = 10949
$ws.Range("U3").Value = -4539
$ws.Range("V3").Value = 141277
$ws.Range("W3").Value = 2.93
$ws.Range("X3").Value = -1.19
$ws.Range("Y3").Value = -2.32
$ws.Range("Z3").Value = -0.86
$ws.Range("AA3").Value = 138.14
$ws.Range("AB3").Value = 10037.45
$ws.Range("AC3").Value = -12164
$ws.Range("AD3").Value = -18.07
$ws.Range("AE3").Value = 549788
$ws.Range("AF3").Value = 0.4
$ws.Range("AG3").Value = 2000
$ws.Range("AH3").Value = 0.91
$ws.Range("AI3").Value = -15.43
$ws.Range("AJ3").Value = 31490892

# Row 4
$ws.Range("D4").Value = 241143
$ws.Range("E4").Value = 7633
$ws.Range("F4").Value = 9404
$ws.Range("G4").Value = 3195
$ws.Range("H4").Value = 2469
$ws.Range("I4").Value = 1682
$ws.Range("J4").Value = 787
$ws.Range("K4").Value = 419159
$ws.Range("L4").Value = 246519
$ws.Range("M4").Value = 172640
$ws.Range("N4").Value = 163045
$ws.Range("O4").Value = 9595
$ws.Range("P4").Value = 1575
$ws.Range("Q4").Value = 15069
$ws.Range("R4").Value = -9684
$ws.Range("S4").Value = -212
$ws.Range("T4").Value = 12070
$ws.Range("U4").Value = 2999
$ws.Range("V4").Value = 145285
$ws.Range("W4").Value = 3.16
$ws.Range("X4").Value = 1.02
$ws.Range("Y4").Value = 1.03
$ws.Range("Z4").Value = 0.6
$ws.Range("AA4").Value = 142.79
$ws.Range("AB4").Value = 10102.2
$ws.Range("AC4").Value = 5341
$ws.Range("AD4").Value = 39.2
$ws.Range("AE4").Value = 551720
$ws.Range("AF4").Value = 0.38
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 0.96
$ws.Range("AI4").Value = 35.14
$ws.Range("AJ4").Value = 31490892

# Row 5
$ws.Range("D5").Value = 179261
$ws.Range("E5").Value = 8010
$ws.Range("F5").Value = 8010
$ws.Range("G5").Value = 2540
$ws.Range("H5").Value = -206
$ws.Range("I5").Value = -1364
$ws.Range("J5").Value = 1158
$ws.Range("K5").Value = 279485
$ws.Range("L5").Value = 145971
$ws.Range("M5").Value = 133514
$ws.Range("N5").Value = 126586
$ws.Range("O5").Value = 6928
$ws.Range("P5").Value = 1406
$ws.Range("Q5").Value = 15693
$ws.Range("R5").Value = -12632
$ws.Range("S5").Value = 749
$ws.Range("T5").Value = 8839
$ws.Range("U5").Value = 6854
$ws.Range("V5").Value = 76778
$ws.Range("W5").Value = 4.47
$ws.Range("X5").Value = -0.12
$ws.Range("Y5").Value = -0.9399999999999999
$ws.Range("Z5").Value = -0.06
$ws.Range("AA5").Value = 109.33
$ws.Range("AB5").Value = 10835.63
$ws.Range("AC5").Value = -4437
$ws.Range("AD5").Value = -44.85
$ws.Range("AE5").Value = 450420
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 5200
$ws.Range("AH5").Value = 2.61
$ws.Range("AI5").Value = -107.15
$ws.Range("AJ5").Value = 28122047

# Row 6
$ws.Range("D6").Value = 178208
$ws.Range("E6").Value = 5970
$ws.Range("F6").Value = 5970
$ws.Range("G6").Value = -939
$ws.Range("H6").Value = -4650
$ws.Range("I6").Value = -5449
$ws.Range("K6").Value = 262593
$ws.Range("L6").Value = 138319
$ws.Range("M6").Value = 124274
$ws.Range("N6").Value = 117054
$ws.Range("P6").Value = 1414
$ws.Range("Q6").Value = 1733
$ws.Range("R6").Value = -2866
$ws.Range("S6").Value = -6926
$ws.Range("T6").Value = 5971
$ws.Range("U6").Value = -4238
$ws.Range("V6").Value = 77956
$ws.Range("W6").Value = 3.35
$ws.Range("X6").Value = -2.61
$ws.Range("Y6").Value = -4.47
$ws.Range("Z6").Value = -1.72
$ws.Range("AA6").Value = 111.3
$ws.Range("AB6").Value = 10297.93
$ws.Range("AC6").Value = -19327
$ws.Range("AD6").Value = -10.92
$ws.Range("AE6").Value = 414048
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 5200
$ws.Range("AH6").Value = 2.46
$ws.Range("AI6").Value = -26.98
$ws.Range("AJ6").Value = 28288755

# Row 7
$ws.Range("D7").Value = 177067
$ws.Range("E7").Value = 5553
$ws.Range("G7").Value = 2566
$ws.Range("H7").Value = 2131
$ws.Range("I7").Value = 1793
$ws.Range("K7").Value = 324398
$ws.Range("L7").Value = 198966
$ws.Range("M7").Value = 125433
$ws.Range("N7").Value = 117141
$ws.Range("P7").Value = 1411
$ws.Range("Q7").Value = 20298
$ws.Range("R7").Value = -3891
$ws.Range("S7").Value = -12710
$ws.Range("T7").Value = 6806
$ws.Range("U7").Value = 9014
$ws.Range("W7").Value = 3.14
$ws.Range("X7").Value = 1.2
$ws.Range("Y7").Value = 1.53
$ws.Range("Z7").Value = 0.73
$ws.Range("AA7").Value = 158.62
$ws.Range("AC7").Value = 6338
$ws.Range("AD7").Value = 18.22
$ws.Range("AE7").Value = 414354
$ws.Range("AF7").Value = 0.28
$ws.Range("AG7").Value = 4680
$ws.Range("AH7").Value = 4.05
$ws.Range("AI7").Value = 73.83

# Row 8
$ws.Range("D8").Value = 179344
$ws.Range("E8").Value = 6405
$ws.Range("G8").Value = 3610
$ws.Range("H8").Value = 2516
$ws.Range("I8").Value = 2420
$ws.Range("K8").Value = 326053
$ws.Range("L8").Value = 199202
$ws.Range("M8").Value = 126850
$ws.Range("N8").Value = 118232
$ws.Range("P8").Value = 1411
$ws.Range("Q8").Value = 16871
$ws.Range("R8").Value = -6825
$ws.Range("S8").Value = -5089
$ws.Range("T8").Value = 6587
$ws.Range("U8").Value = 8042
$ws.Range("W8").Value = 3.57
$ws.Range("X8").Value = 1.4
$ws.Range("Y8").Value = 2.06
$ws.Range("Z8").Value = 0.77
$ws.Range("AA8").Value = 157.04
$ws.Range("AC8").Value = 8556
$ws.Range("AD8").Value = 13.5
$ws.Range("AE8").Value = 418214
$ws.Range("AF8").Value = 0.28
$ws.Range("AG8").Value = 4550
$ws.Range("AH8").Value = 3.94
$ws.Range("AI8").Value = 53.18

# Row 9
$ws.Range("D9").Value = 182561
$ws.Range("E9").Value = 6982
$ws.Range("G9").Value = 4323
$ws.Range("H9").Value = 3171
$ws.Range("I9").Value = 2765
$ws.Range("K9").Value = 329318
$ws.Range("L9").Value = 200528
$ws.Range("M9").Value = 128788
$ws.Range("N9").Value = 119773
$ws.Range("P9").Value = 1411
$ws.Range("Q9").Value = 16382
$ws.Range("R9").Value = -6288
$ws.Range("S9").Value = -5052
$ws.Range("T9").Value = 6219
$ws.Range("U9").Value = 8423
$ws.Range("W9").Value = 3.82
$ws.Range("X9").Value = 1.74
$ws.Range("Y9").Value = 2.32
$ws.Range("Z9").Value = 0.97
$ws.Range("AA9").Value = 155.7
$ws.Range("AC9").Value = 9773
$ws.Range("AD9").Value = 11.82
$ws.Range("AE9").Value = 423664
$ws.Range("AF9").Value = 0.27
$ws.Range("AG9").Value = 4570
$ws.Range("AH9").Value = 3.96
$ws.Range("AI9").Value = 46.76

